$d = $word.ActiveDocument

# The closing/signature block ("Semarang, {tanggalPembuatan}" through
# "{namaPemilik}") needs a left indent of 4320 twips (= 216 pt = 3")
# applied to each paragraph so the field lines up under the tab stop.
$startIndex = 17
$endIndex = 25

for ($i = $startIndex; $i -le $endIndex; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Format.LeftIndent = 216
}
